$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20-26 no longer exist in the updated sheet (dimension shrinks to A1:F19)
$ws.Rows("20:26").Delete()

# Row 2
$ws.Range("B2").Value = "NSE:ABB"
$ws.Range("C2").Value = "NSE:AARTIDRUGS"
$ws.Range("E2").Value = "NSE:MOTHERSON"
$ws.Range("F2").Value = "NSE:GRASIM"

# Row 3
$ws.Range("B3").Value = "NSE:ANANTRAJ"
$ws.Range("C3").Value = "NSE:ADROITINFO"
$ws.Range("F3").Value = "NSE:HDFCAMC"

# Row 4
$ws.Range("B4").Value = "NSE:GRASIM"
$ws.Range("C4").Value = "NSE:AFFLE"
$ws.Range("F4").Value = "NSE:PFC"

# Row 5
$ws.Range("B5").Value = "NSE:HAL"
$ws.Range("C5").Value = "NSE:ANMOL"

# Row 6
$ws.Range("B6").Value = "NSE:HDFCAMC"
$ws.Range("C6").Value = "NSE:AXISBANK"

# Row 7
$ws.Range("B7").Value = "NSE:INGERRAND"
$ws.Range("C7").Value = "NSE:BANDHANBNK"

# Row 8
$ws.Range("B8").Value = "NSE:IONEXCHANG"
$ws.Range("C8").Value = "NSE:BBOX"

# Row 9
$ws.Range("B9").Value = "NSE:ISEC"
$ws.Range("C9").Value = "NSE:EMMBI"

# Row 10
$ws.Range("B10").Value = "NSE:JKIL"
$ws.Range("C10").Value = "NSE:GTLINFRA"

# Row 11
$ws.Range("B11").Value = "NSE:LTTS"
$ws.Range("C11").Value = "NSE:HAVISHA"

# Row 12
$ws.Range("B12").Value = "NSE:MATRIMONY"
$ws.Range("C12").Value = "NSE:HEUBACHIND"

# Row 13
$ws.Range("B13").Value = "NSE:MONARCH"
$ws.Range("C13").Value = "NSE:HUHTAMAKI"

# Row 14
$ws.Range("B14").Value = "NSE:MOREPENLAB"
$ws.Range("C14").Value = "NSE:JAGSNPHARM"

# Row 15
$ws.Range("B15").Value = "NSE:RADIANTCMS"
$ws.Range("C15").Value = "NSE:LOKESHMACH"

# Row 16
$ws.Range("B16").Value = "NSE:RATNAMANI"
$ws.Range("C16").Value = "NSE:MAZDA"

# Row 17
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = "NSE:NUCLEUS"

# Row 18
$ws.Range("B18").Value = $null
$ws.Range("C18").Value = "NSE:OLECTRA"

# Row 19
$ws.Range("B19").Value = $null
$ws.Range("C19").Value = "NSE:QUICKHEAL"
